$p = $ppt.ActivePresentation
$m = $p.SlideMaster
Write-Host "before:" $m.Name
$m.Name = "MyMaster"
Write-Host "after:" $m.Name
